$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 71.80040646296555
$ws.Range("C2").Value = 76.0515352652897
$ws.Range("D2").Value = 67.08205777666116
$ws.Range("E2").Value = 81.3979961982961

$ws.Range("B3").Value = 94.79333533292099
$ws.Range("C3").Value = 92.92205992593887
$ws.Range("D3").Value = 94.23739774346879
$ws.Range("E3").Value = 94.1575531821819

$ws.Range("B4").Value = 99.02321707674155
$ws.Range("C4").Value = 98.81177698375096
$ws.Range("D4").Value = 98.99818792400407
$ws.Range("E4").Value = 99.07712320080182

$ws.Range("B5").Value = 98.82435018620338
$ws.Range("C5").Value = 98.94723411352066
$ws.Range("D5").Value = 98.81146515271119
$ws.Range("E5").Value = 98.51576391116924

$ws.Range("B6").Value = 98.53849430365493
$ws.Range("C6").Value = 98.47104428594164
$ws.Range("D6").Value = 98.47552954347444
$ws.Range("E6").Value = 98.42869838002257

$ws.Range("B7").Value = 97.52012708421618
$ws.Range("C7").Value = 97.52601928084786
$ws.Range("D7").Value = 97.56723329164841
$ws.Range("E7").Value = 97.57435071890701

$ws.Range("B8").Value = 96.17949301444447
$ws.Range("C8").Value = 96.14312341725785
$ws.Range("D8").Value = 96.14193776916262
$ws.Range("E8").Value = 96.10174774157602
